$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny precision change in A2 (date/time serial value)
$ws.Range("A2").Value = 45875.04187768519

# Add new row 3 with data
$ws.Range("A3").Value = 45875.08356718349
$ws.Range("B3").Value = 2025
$ws.Range("C3").Value = 23
$ws.Range("D3").Value = 13.72
$ws.Range("E3").Value = 92.37
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 4.12
$ws.Range("H3").Value = "E"
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "02:00:20"

# Copy style from A2 to A3 so date formatting matches
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
